# Updated cryptos list (Price / Volume(1h) columns) as scraped by the
# GitHub Actions job. Prices that look numeric are forced back to Text
# (NumberFormat "@") before assignment so they round-trip the same way
# the source data does (e.g. "589.42" stays a string, not a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.725.26"
$ws.Range("E2").Value = "  +3.70%  "

$ws.Range("D3").Value = "3.686.74"
$ws.Range("E3").Value = "  +8.27%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.42"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.05"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("D7").Value = "3.677.88"
$ws.Range("E7").Value = "  +8.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  +4.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.201"
$ws.Range("E10").Value = "  +1.54%  "

$ws.Range("E11").Value = "  +4.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.11"
$ws.Range("E12").Value = "  +3.56%  "

$ws.Range("E13").Value = "  +2.37%  "

$ws.Range("D14").Value = "4.283.54"
$ws.Range("E14").Value = "  +8.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "683.80"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.05"
$ws.Range("E16").Value = "  +5.08%  "

$ws.Range("D17").Value = "3.690.50"
$ws.Range("E17").Value = "  +8.55%  "

$ws.Range("D18").Value = "71.830.66"
$ws.Range("E18").Value = "  +3.69%  "

$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.18"
$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("E21").Value = "  +3.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.946"
$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.34"
$ws.Range("E23").Value = "  +17.68%  "

$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "104.10"
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("E26").Value = "  +4.05%  "

$ws.Range("E27").Value = "  +6.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  +5.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.45"
$ws.Range("E29").Value = "  +5.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.28"
$ws.Range("E30").Value = "  +6.17%  "

$ws.Range("E31").Value = "  +7.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.24"
$ws.Range("E32").Value = "  +13.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.34"
$ws.Range("E33").Value = "  +3.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "566.77"
$ws.Range("E34").Value = "  +2.00%  "

$ws.Range("E35").Value = "  +4.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.52"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("D37").Value = "3.805.83"
$ws.Range("E37").Value = "  +5.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  +5.11%  "

$ws.Range("D40").Value = "0.0₃0778"
$ws.Range("E40").Value = "  +5.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.65"
$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("E42").Value = "  +6.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0469"
$ws.Range("E43").Value = "  +10.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  +4.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.354"
$ws.Range("E45").Value = "  +5.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  +8.96%  "

$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("E48").Value = "  +4.29%  "

$ws.Range("E49").Value = "  +3.23%  "

$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.83"
$ws.Range("E51").Value = "  +2.88%  "
